# Updates cryptos list prices (column D) and 1h volume deltas (column E).
# D-column values are digit/period strings that Excel would otherwise
# auto-convert to a Double on assignment (losing the exact text and the
# ungrouped-cell style), so they are written with a leading apostrophe to
# force text, then the quote-prefix style introduced by that is reset back
# to "Normal" so the cell keeps its original (unstyled) appearance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.782.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.69%  "

$ws.Range("D3").Value = "'3.324.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.83%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'604.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.66%  "

$ws.Range("D6").Value = "'142.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.31%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").Value = "'3.322.41"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.98%  "

$ws.Range("D9").Value = "'0.520"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.25%  "

$ws.Range("E10").Value = "  +1.51%  "

$ws.Range("D11").Value = "'5.53"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.14%  "

$ws.Range("D12").Value = "'0.470"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.10%  "

$ws.Range("E13").Value = "  -0.33%  "

$ws.Range("D14").Value = "'35.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.43%  "

$ws.Range("D15").Value = "'3.871.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.24%  "

$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("D17").Value = "'3.323.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.98%  "

$ws.Range("D18").Value = "'63.851.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.96%  "

$ws.Range("D19").Value = "'6.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.52%  "

$ws.Range("D20").Value = "'480.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.34%  "

$ws.Range("D21").Value = "'14.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("E22").Value = "  +3.31%  "

$ws.Range("E23").Value = "  +1.34%  "

$ws.Range("E24").Value = "  +4.12%  "

$ws.Range("D25").Value = "'84.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.78%  "

$ws.Range("E27").Value = "  +1.42%  "

$ws.Range("D28").Value = "'8.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.78%  "

$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("E30").Value = "  -0.83%  "

$ws.Range("E31").Value = "  +2.09%  "

$ws.Range("D32").Value = "'28.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.53%  "

$ws.Range("E33").Value = "  -1.64%  "

$ws.Range("D34").Value = "'2.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.16%  "

$ws.Range("E35").Value = "  +1.39%  "

$ws.Range("D36").Value = "'6.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.69%  "

$ws.Range("D37").Value = "'52.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.89%  "

$ws.Range("E38").Value = "  +3.41%  "

$ws.Range("E39").Value = "  +2.21%  "

$ws.Range("D40").Value = "'434.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.37%  "

$ws.Range("D41").Value = "'3.096.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.87%  "

$ws.Range("E42").Value = "  +6.41%  "

$ws.Range("D43").Value = "'2.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.51%  "

$ws.Range("D44").Value = "'8.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.81%  "

$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").Value = "'2.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.01%  "

$ws.Range("D47").Value = "'37.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +14.00%  "

$ws.Range("D48").Value = "'26.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.15%  "

$ws.Range("E50").Value = "  +1.20%  "

$ws.Range("E51").Value = "  -0.28%  "
